$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Widen column D from 28 to 32 characters
$ws.Columns.Item(4).ColumnWidth = 31.1

# 2. Insert the new row 5 ("生成AI案件を回してくれるパートナー募集!")
$ws.Rows.Item(5).Insert()

# 3. Insert two new rows at 13 and 14 (PHP request, then Stripe duplicate)
#    old row 12 (Stripe/5493449) is pushed down to row 15 automatically
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(14).Insert()

# 4. Refresh the "取得日時" timestamp on every existing data row (2-12 pre-shift == now 2,3,4,6-12,15)

# Row 2 (unchanged content, timestamp refresh)
$ws.Cells.Item(2,1).Value = '2026-02-17 12:58:45'

# Row 3 (unchanged content, timestamp refresh)
$ws.Cells.Item(3,1).Value = '2026-02-17 12:58:45'

# Row 4 (unchanged content, timestamp refresh)
$ws.Cells.Item(4,1).Value = '2026-02-17 12:58:45'

# Row 5 (brand new posting - 生成AI案件)
$ws.Cells.Item(5,1).Value = '2026-02-17 12:58:45'
$ws.Cells.Item(5,2).Value = '【急募】生成AI案件を回してくれるパートナー募集!'
$ws.Cells.Item(5,3).Value = 'システム開発'
$ws.Cells.Item(5,4).Value = '1,000,000 円 ~ 3,000,000 円 / 固定'
$ws.Cells.Item(5,5).Value = '期限情報なし'
$ws.Cells.Item(5,6).Value = 'https://www.lancers.jp/work/detail/5493776'
$ws.Cells.Item(5,7).Value = 310
$ws.Cells.Item(5,8).Value = '🔥AI,Ai'

# Row 6 (was row 5, timestamp refresh)
$ws.Cells.Item(6,1).Value = '2026-02-17 12:58:45'

# Row 7 (was row 6, timestamp refresh)
$ws.Cells.Item(7,1).Value = '2026-02-17 12:58:45'

# Row 8 (was row 7, timestamp refresh)
$ws.Cells.Item(8,1).Value = '2026-02-17 12:58:45'

# Row 9 (was row 8, timestamp refresh)
$ws.Cells.Item(9,1).Value = '2026-02-17 12:58:45'

# Row 10 (was row 9, timestamp refresh)
$ws.Cells.Item(10,1).Value = '2026-02-17 12:58:45'

# Row 11 (was row 10, timestamp refresh)
$ws.Cells.Item(11,1).Value = '2026-02-17 12:58:45'

# Row 12 (was row 11, timestamp refresh)
$ws.Cells.Item(12,1).Value = '2026-02-17 12:58:45'

# Row 13 (brand new posting - PHP verification)
$ws.Cells.Item(13,1).Value = '2026-02-17 12:58:45'
$ws.Cells.Item(13,2).Value = '【急募】PHPバージョンアップ検証のためのテスト環境構築依頼'
$ws.Cells.Item(13,3).Value = 'システム開発'
$ws.Cells.Item(13,4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(13,5).Value = '期限情報なし'
$ws.Cells.Item(13,6).Value = 'https://www.lancers.jp/work/detail/5493555'
$ws.Cells.Item(13,7).Value = 28
$ws.Cells.Item(13,8).Value = '○PHP'

# Row 14 (brand new posting - Stripe duplicate, url 5493650)
$ws.Cells.Item(14,1).Value = '2026-02-17 12:58:45'
$ws.Cells.Item(14,2).Value = '【設計済み!作業時間~10時間】Stripe(銀行振込)を用いた月額課金システムの構築'
$ws.Cells.Item(14,3).Value = 'システム開発'
$ws.Cells.Item(14,4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(14,5).Value = '期限情報なし'
$ws.Cells.Item(14,6).Value = 'https://www.lancers.jp/work/detail/5493650'
$ws.Cells.Item(14,7).Value = 28

# Row 15 (was row 12, timestamp refresh)
$ws.Cells.Item(15,1).Value = '2026-02-17 12:58:45'

# Row 16 (brand new posting - AWS domain integration, appended at end)
$ws.Cells.Item(16,1).Value = '2026-02-17 12:58:45'
$ws.Cells.Item(16,2).Value = '初回 【自社HP構築】AWSを活かした最適なドメイン統合構成の設計・実装支援'
$ws.Cells.Item(16,3).Value = 'システム開発'
$ws.Cells.Item(16,4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(16,5).Value = '期限情報なし'
$ws.Cells.Item(16,6).Value = 'https://www.lancers.jp/work/detail/5493714'
$ws.Cells.Item(16,7).Value = 18

# 5. Rebuild hyperlinks for F2:F16 (style + relationship) in row order
$ws.Hyperlinks.Delete()
for ($r = 2; $r -le 16; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($cell, $cell.Value)
    $cell.Style = "Hyperlink"
}

